# [CAN2-16] TLD Reqs Rev1 ready for review per [CAN2-24]
# Fix the "instansiate" -> "instantiate" typo in the Requirement Body column
# (D2:D9, the Initialization / TOP_INIT_xx requirements) and leave the
# worksheet with D3 as the active/selected cell, matching the reviewer's
# final click-through.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$fixes = @{
    "D2" = "The module shall instantiate the Microcontroller Bus Interface as DUT0."
    "D3" = "The module shall instantiate the Configuration Registers as DUT1."
    "D4" = "The module shall instantiate the Tx FIFO as DUT2."
    "D5" = "The module shall instantiate the Rx FIFO as DUT3."
    "D6" = "The module shall instantiate the Tx Priority Logic as DUT4."
    "D7" = "The module shall instantiate the Acceptance Filter as DUT5."
    "D8" = "The module shall instantiate the Bit Stream Processor as DUT6."
    "D9" = "The module shall instantiate the Bit Timing Module as DUT7."
}

foreach ($addr in $fixes.Keys) {
    $ws.Range($addr).Value = $fixes[$addr]
}

# Leave the selection on D3, same as the final state captured in the workbook.
$ws.Activate() | Out-Null
$ws.Range("D3").Select() | Out-Null
